$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C81").Value = "Calmar_Ratio_test1"
$ws.Range("A81").Value = "Calmar Ratio1"
$ws.Range("B81").Value = "Test Calmar ratio with scale=1"

$ws.Range("A82").Value = "Calmar Ratio2"
$ws.Range("B82").Value = "Test Calmar ratio with scale=252"
$ws.Range("C82").Value = "Calmar_Ratio_test2"

$ws.Range("E86").Select()
